# Auto-generated edit script: update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell holding the default (unstyled) format, used to restore
# style after temporarily marking numeric-looking text cells as Text
# so Excel does not auto-convert them into real numbers.
$donorStyle = $ws.Range("D36").Style

$ws.Range("D2").Value = "36.978.34"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "2.058.32"
$ws.Range("E3").Value = "  -1.95%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.03"
$ws.Range("D5").Style = $donorStyle
$ws.Range("E5").Value = "  -1.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.671"
$ws.Range("D6").Style = $donorStyle
$ws.Range("E6").Value = "  +1.33%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.36"
$ws.Range("D8").Style = $donorStyle
$ws.Range("E8").Value = "  +12.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.49"
$ws.Range("D9").Style = $donorStyle
$ws.Range("E9").Value = "  -0.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.381"
$ws.Range("D10").Style = $donorStyle
$ws.Range("E10").Value = "  +0.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("D11").Style = $donorStyle
$ws.Range("E11").Value = "  +8.02%  "

$ws.Range("E12").Value = "  -3.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.97"
$ws.Range("D13").Style = $donorStyle
$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("D14").Value = "2.362.16"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.814"
$ws.Range("D15").Style = $donorStyle
$ws.Range("E15").Value = "  -3.13%  "

$ws.Range("E16").Value = "  +2.34%  "

$ws.Range("D17").Value = "2.055.68"
$ws.Range("E17").Value = "  -2.70%  "

$ws.Range("D18").Value = "36.911.31"
$ws.Range("E18").Value = "  +0.47%  "

$ws.Range("D19").Value = "0.0₃0938"
$ws.Range("E19").Value = "  +13.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.65"
$ws.Range("D20").Style = $donorStyle
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.25"
$ws.Range("D21").Style = $donorStyle
$ws.Range("E21").Value = "  +6.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.37"
$ws.Range("D22").Style = $donorStyle
$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.42"
$ws.Range("D23").Style = $donorStyle

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = $donorStyle
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("E25").Value = "  -5.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.56"
$ws.Range("D26").Style = $donorStyle
$ws.Range("E26").Value = "  -0.75%  "

$ws.Range("E27").Value = "  -3.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.09"
$ws.Range("D28").Style = $donorStyle
$ws.Range("E28").Value = "  -5.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.99"
$ws.Range("D29").Style = $donorStyle
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("E30").Value = "  +1.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.58"
$ws.Range("D31").Style = $donorStyle
$ws.Range("E31").Value = "  +1.35%  "

$ws.Range("E32").Value = "  +5.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0628"
$ws.Range("D33").Style = $donorStyle
$ws.Range("E33").Value = "  +1.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.37"
$ws.Range("D34").Style = $donorStyle
$ws.Range("E34").Value = "  +6.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0879"
$ws.Range("D35").Style = $donorStyle
$ws.Range("E35").Value = "  -5.06%  "

$ws.Range("E37").Value = "  -6.05%  "

$ws.Range("E38").Value = "  -4.18%  "

$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("E40").Value = "  +24.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.94"
$ws.Range("D41").Style = $donorStyle
$ws.Range("E41").Value = "  +6.85%  "

$ws.Range("E42").Value = "  -0.39%  "

$ws.Range("E43").Value = "  -2.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.84"
$ws.Range("D44").Style = $donorStyle
$ws.Range("E44").Value = "  -1.55%  "

$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.12"
$ws.Range("D46").Style = $donorStyle
$ws.Range("E46").Value = "  +42.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "14.02"
$ws.Range("D47").Style = $donorStyle
$ws.Range("E47").Value = "  -47.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.34"
$ws.Range("D48").Style = $donorStyle
$ws.Range("E48").Value = "  +10.66%  "

$ws.Range("E49").Value = "  +6.59%  "

$ws.Range("D50").Value = "1.300.17"
$ws.Range("E50").Value = "  -3.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.92"
$ws.Range("D51").Style = $donorStyle
$ws.Range("E51").Value = "  +0.72%  "

